$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.9230769230769231
$ws.Range("D4").Value = 0.9487179487179487
